$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4115256666666666
$ws.Range("H2").Value = 1.234577
$ws.Range("I2").Value = 0.2245998342667577
$ws.Range("J2").Value = 0.2245998342667577
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.9721959999999999
$ws.Range("N2").Value = 2.916588
$ws.Range("O2").Value = 0.7027023771175303
$ws.Range("P2").Value = 0.7027023771175303
$ws.Range("Q2").Value = 0.4000836070306666
$ws.Range("R2").Value = 3.600752463275999
$ws.Range("S2").Value = 0.157826837439454
$ws.Range("T2").Value = 0.157826837439454
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4115256666666666
$ws.Range("H3").Value = 1.234577
$ws.Range("I3").Value = 0.2245998342667577
$ws.Range("J3").Value = 0.2245998342667577
$ws.Range("O3").Value = 0.1592492623233027
$ws.Range("P3").Value = 0.1592492623233027
$ws.Range("Q3").Value = 0.09066856945699997
$ws.Range("R3").Value = 0.8160171251129997
$ws.Range("S3").Value = 0.03576735792491722
$ws.Range("T3").Value = 0.03576735792491722
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4115256666666666
$ws.Range("H4").Value = 1.234577
$ws.Range("I4").Value = 0.2245998342667577
$ws.Range("J4").Value = 0.2245998342667577
$ws.Range("M4").Value = 0.1909913333333333
$ws.Range("N4").Value = 0.572974
$ws.Range("O4").Value = 0.138048360559167
$ws.Range("P4").Value = 0.138048360559167
$ws.Range("Q4").Value = 0.07859783577755554
$ws.Range("R4").Value = 0.7073805219979998
$ws.Range("S4").Value = 0.03100563890238652
$ws.Range("T4").Value = 0.03100563890238652
$ws.Range("I5").Value = 0.3944722233087159
$ws.Range("J5").Value = 0.3944722233087159
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.9721959999999999
$ws.Range("N5").Value = 2.916588
$ws.Range("O5").Value = 0.7027023771175303
$ws.Range("P5").Value = 0.7027023771175303
$ws.Range("Q5").Value = 0.7026802601613333
$ws.Range("R5").Value = 6.324122341452
$ws.Range("S5").Value = 0.2771965690258719
$ws.Range("T5").Value = 0.2771965690258719
$ws.Range("I6").Value = 0.3944722233087159
$ws.Range("J6").Value = 0.3944722233087159
$ws.Range("O6").Value = 0.1592492623233027
$ws.Range("P6").Value = 0.1592492623233027
$ws.Range("S6").Value = 0.06281941056894615
$ws.Range("T6").Value = 0.06281941056894615
$ws.Range("I7").Value = 0.3944722233087159
$ws.Range("J7").Value = 0.3944722233087159
$ws.Range("M7").Value = 0.1909913333333333
$ws.Range("N7").Value = 0.572974
$ws.Range("O7").Value = 0.138048360559167
$ws.Range("P7").Value = 0.138048360559167
$ws.Range("Q7").Value = 0.1380440156051111
$ws.Range("R7").Value = 1.242396140446
$ws.Range("S7").Value = 0.05445624371389785
$ws.Range("T7").Value = 0.05445624371389785
$ws.Range("G8").Value = 0.6979596666666668
$ws.Range("H8").Value = 2.093879
$ws.Range("I8").Value = 0.3809279424245264
$ws.Range("J8").Value = 0.3809279424245264
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.9721959999999999
$ws.Range("N8").Value = 2.916588
$ws.Range("O8").Value = 0.7027023771175303
$ws.Range("P8").Value = 0.7027023771175303
$ws.Range("Q8").Value = 0.6785535960946667
$ws.Range("R8").Value = 6.106982364852001
$ws.Range("S8").Value = 0.2676789706522044
$ws.Range("T8").Value = 0.2676789706522044
$ws.Range("G9").Value = 0.6979596666666668
$ws.Range("H9").Value = 2.093879
$ws.Range("I9").Value = 0.3809279424245264
$ws.Range("J9").Value = 0.3809279424245264
$ws.Range("O9").Value = 0.1592492623233027
$ws.Range("P9").Value = 0.1592492623233027
$ws.Range("Q9").Value = 0.153776567639
$ws.Range("R9").Value = 1.383989108751
$ws.Range("S9").Value = 0.06066249382943936
$ws.Range("T9").Value = 0.06066249382943936
$ws.Range("G10").Value = 0.6979596666666668
$ws.Range("H10").Value = 2.093879
$ws.Range("I10").Value = 0.3809279424245264
$ws.Range("J10").Value = 0.3809279424245264
$ws.Range("M10").Value = 0.1909913333333333
$ws.Range("N10").Value = 0.572974
$ws.Range("O10").Value = 0.138048360559167
$ws.Range("P10").Value = 0.138048360559167
$ws.Range("Q10").Value = 0.1333042473495556
$ws.Range("R10").Value = 1.199738226146
$ws.Range("S10").Value = 0.05258647794288263
$ws.Range("T10").Value = 0.05258647794288263
